$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 holds three rotating "blog" slots (ser numbers). A new blog post (115)
# was published, so each slot's series number is bumped by one: the oldest
# (112) is dropped and 115 becomes the newest entry.
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 113"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 114"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 115"
